$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text (avoids Excel auto-converting
# numeric-looking strings like "610.17" or "0.0961" into floating point
# numbers), while keeping the cell free of any extra style reference.
function Set-TextCell {
    param($rng, $val)
    $rng.Value = "'" + $val
    $rng.Style = "Normal"
}

# Row-level updates: cell letter -> new value, keyed by row number
$rowUpdates = @{
    2 = @{ D = "63.368.65"; E = "  +0.23%  " }
    3 = @{ D = "2.658.49"; E = "  +3.25%  " }
    4 = @{ E = "  +0.03%  " }
    5 = @{ D = "610.17"; E = "  +3.99%  " }
    6 = @{ D = "143.55"; E = "  -0.75%  " }
    7 = @{ E = "  +0.02%  " }
    8 = @{ D = "0.586"; E = "  -0.75%  " }
    9 = @{ D = "2.657.67"; E = "  +3.29%  " }
    10 = @{ E = "  +0.56%  " }
    11 = @{ D = "5.62"; E = "  -0.10%  " }
    12 = @{ E = "  +0.31%  " }
    13 = @{ D = "0.362"; E = "  +3.11%  " }
    14 = @{ D = "27.29"; E = "  +0.40%  " }
    15 = @{ D = "3.128.01"; E = "  +3.02%  " }
    16 = @{ D = "63.211.11"; E = "  +0.11%  " }
    17 = @{ E = "  -0.99%  " }
    18 = @{ D = "2.661.69"; E = "  +3.05%  " }
    19 = @{ D = "11.45"; E = "  +3.41%  " }
    20 = @{ D = "342.00"; E = "  +0.18%  " }
    21 = @{ D = "4.42"; E = "  +1.83%  " }
    22 = @{ D = "6.87"; E = "  +3.36%  " }
    23 = @{ D = "1.00"; E = "  -0.08%  " }
    24 = @{ D = "66.94"; E = "  -1.46%  " }
    25 = @{ E = "  +2.13%  " }
    26 = @{ D = "1.55"; E = "  -0.52%  " }
    27 = @{ D = "8.64"; E = "  +5.08%  " }
    28 = @{ E = "  -0.52%  " }
    29 = @{ D = "547.27"; E = "  +16.19%  " }
    30 = @{ D = "0.998"; E = "  -0.27%  " }
    31 = @{ D = "7.82"; E = "  -1.52%  " }
    32 = @{ D = "2.06"; E = "  +6.09%  " }
    33 = @{ E = "  +6.88%  " }
    34 = @{ D = "0.0₃0807"; E = "  +0.68%  " }
    35 = @{ D = "172.94"; E = "  -1.76%  " }
    36 = @{ D = "5.16"; E = "  +13.67%  " }
    37 = @{ D = "0.407"; E = "  +1.58%  " }
    38 = @{ D = "1.00"; E = "  -0.08%  " }
    39 = @{ D = "19.15"; E = "  +1.48%  " }
    40 = @{ E = "  +8.68%  " }
    41 = @{ D = "175.06"; E = "  +11.01%  " }
    42 = @{ D = "0.997"; E = "  -0.19%  " }
    43 = @{ D = "3.74"; E = "  +1.44%  " }
    44 = @{ D = "22.20"; E = "  +4.39%  " }
    45 = @{ E = "  +6.18%  " }
    46 = @{ D = "0.631"; E = "  -0.34%  " }
    47 = @{ B = "Stellar"; C = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D = "0.0961"; E = "  -0.20%  " }
    48 = @{ B = "VeChain"; C = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D = "0.0240"; E = "  +1.26%  " }
    49 = @{ D = "18.70"; E = "  +2.92%  " }
    50 = @{ D = "1.75"; E = "  +4.44%  " }
    51 = @{ D = "11.31"; E = "  -0.55%  " }
}

foreach ($rowNum in $rowUpdates.Keys) {
    $cellVals = $rowUpdates[$rowNum]
    foreach ($col in $cellVals.Keys) {
        Set-TextCell $ws.Range("$col$rowNum") $cellVals[$col]
    }
}
